{"js": "const replacements = [\n  { row: 0, col: 0, oldText: \"87\u00f75=17, 2\", newText: \"50\u00f72=25, 0\" },\n  { row: 0, col: 1, oldText: \"22\u00f73=7, 1\", newText: \"45\u00f79=5, 0\" },\n  { row: 0, col: 2, oldText: \"25\u00f73=8, 1\", newText: \"36\u00f74=9, 0\" },\n  { row: 0, col: 3, oldText: \"54\u00f74=13, 2\", newText: \"74\u00f77=10, 4\" },\n  { row: 0, col: 4, oldText: \"30\u00f75=6, 0\", newText: \"63\u00f75=12, 3\" },\n  { row: 4, col: 0, oldText: \"43\u00f76=7, 1\", newText: \"40\u00f78=5, 0\" },\n  { row: 4, col: 1, oldText: \"77\u00f76=12, 5\", newText: \"49\u00f72=24, 1\" },\n  { row: 4, col: 2, oldText: \"51\u00f72=25, 1\", newText: \"18\u00f77=2, 4\" },\n  { row: 4, col: 3, oldText: \"43\u00f72=21, 1\", newText: \"35\u00f78=4, 3\" },\n  { row: 4, col: 4, oldText: \"60\u00f74=15, 0\", newText: \"81\u00f75=16, 1\" },\n  { row: 8, col: 0, oldText: \"52\u00f73=17, 1\", newText: \"60\u00f73=20, 0\" },\n  { row: 8, col: 1, oldText: \"18\u00f77=2, 4\", newText: \"29\u00f74=7, 1\" },\n  { row: 8, col: 2, oldText: \"62\u00f78=7, 6\", newText: \"86\u00f77=12, 2\" },\n  { row: 8, col: 3, oldText: \"55\u00f72=27, 1\", newText: \"56\u00f72=28, 0\" },\n  { row: 8, col: 4, oldText: \"29\u00f78=3, 5\", newText: \"40\u00f76=6, 4\" },\n  { row: 12, col: 0, oldText: \"38\u00f77=5, 3\", newText: \"23\u00f72=11, 1\" },\n  { row: 12, col: 1, oldText: \"91\u00f76=15, 1\", newText: \"23\u00f74=5, 3\" },\n  { row: 12, col: 2, oldText: \"89\u00f72=44, 1\", newText: \"76\u00f77=10, 6\" },\n  { row: 12, col: 3, oldText: \"69\u00f73=23, 0\", newText: \"24\u00f73=8, 0\" },\n  { row: 12, col: 4, oldText: \"35\u00f75=7, 0\", newText: \"57\u00f75=11, 2\" },\n  { row: 16, col: 0, oldText: \"70\u00f78=8, 6\", newText: \"39\u00f78=4, 7\" },\n  { row: 16, col: 1, oldText: \"28\u00f74=7, 0\", newText: \"42\u00f76=7, 0\" },\n  { row: 16, col: 2, oldText: \"56\u00f78=7, 0\", newText: \"68\u00f79=7, 5\" },\n  { row: 16, col: 3, oldText: \"63\u00f79=7, 0\", newText: \"48\u00f73=16, 0\" },\n  { row: 16, col: 4, oldText: \"85\u00f76=14, 1\", newText: \"99\u00f72=49, 1\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  // Scope the search to this single cell's body so that a newly written\n  // value can never be re-matched by a later replacement rule (some of the\n  // old/new strings in this sheet coincide across different cells).\n  const searchResults = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match of \"${oldText}\" in cell (${row}, ${col}), found ${searchResults.items.length}.`\n    );\n  }\n\n  searchResults.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# (row, col) use 1-based indices as in the Word object model; oldText is asserted\n# before writing to guard against any cell/table drift.\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"87\u00f75=17, 2\"; NewText = \"50\u00f72=25, 0\" }\n    @{ Row = 1; Col = 2; OldText = \"22\u00f73=7, 1\"; NewText = \"45\u00f79=5, 0\" }\n    @{ Row = 1; Col = 3; OldText = \"25\u00f73=8, 1\"; NewText = \"36\u00f74=9, 0\" }\n    @{ Row = 1; Col = 4; OldText = \"54\u00f74=13, 2\"; NewText = \"74\u00f77=10, 4\" }\n    @{ Row = 1; Col = 5; OldText = \"30\u00f75=6, 0\"; NewText = \"63\u00f75=12, 3\" }\n    @{ Row = 5; Col = 1; OldText = \"43\u00f76=7, 1\"; NewText = \"40\u00f78=5, 0\" }\n    @{ Row = 5; Col = 2; OldText = \"77\u00f76=12, 5\"; NewText = \"49\u00f72=24, 1\" }\n    @{ Row = 5; Col = 3; OldText = \"51\u00f72=25, 1\"; NewText = \"18\u00f77=2, 4\" }\n    @{ Row = 5; Col = 4; OldText = \"43\u00f72=21, 1\"; NewText = \"35\u00f78=4, 3\" }\n    @{ Row = 5; Col = 5; OldText = \"60\u00f74=15, 0\"; NewText = \"81\u00f75=16, 1\" }\n    @{ Row = 9; Col = 1; OldText = \"52\u00f73=17, 1\"; NewText = \"60\u00f73=20, 0\" }\n    @{ Row = 9; Col = 2; OldText = \"18\u00f77=2, 4\"; NewText = \"29\u00f74=7, 1\" }\n    @{ Row = 9; Col = 3; OldText = \"62\u00f78=7, 6\"; NewText = \"86\u00f77=12, 2\" }\n    @{ Row = 9; Col = 4; OldText = \"55\u00f72=27, 1\"; NewText = \"56\u00f72=28, 0\" }\n    @{ Row = 9; Col = 5; OldText = \"29\u00f78=3, 5\"; NewText = \"40\u00f76=6, 4\" }\n    @{ Row = 13; Col = 1; OldText = \"38\u00f77=5, 3\"; NewText = \"23\u00f72=11, 1\" }\n    @{ Row = 13; Col = 2; OldText = \"91\u00f76=15, 1\"; NewText = \"23\u00f74=5, 3\" }\n    @{ Row = 13; Col = 3; OldText = \"89\u00f72=44, 1\"; NewText = \"76\u00f77=10, 6\" }\n    @{ Row = 13; Col = 4; OldText = \"69\u00f73=23, 0\"; NewText = \"24\u00f73=8, 0\" }\n    @{ Row = 13; Col = 5; OldText = \"35\u00f75=7, 0\"; NewText = \"57\u00f75=11, 2\" }\n    @{ Row = 17; Col = 1; OldText = \"70\u00f78=8, 6\"; NewText = \"39\u00f78=4, 7\" }\n    @{ Row = 17; Col = 2; OldText = \"28\u00f74=7, 0\"; NewText = \"42\u00f76=7, 0\" }\n    @{ Row = 17; Col = 3; OldText = \"56\u00f78=7, 0\"; NewText = \"68\u00f79=7, 5\" }\n    @{ Row = 17; Col = 4; OldText = \"63\u00f79=7, 0\"; NewText = \"48\u00f73=16, 0\" }\n    @{ Row = 17; Col = 5; OldText = \"85\u00f76=14, 1\"; NewText = \"99\u00f72=49, 1\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $table.Cell($item.Row, $item.Col)\n    $cellRange = $cell.Range\n    # Trim the trailing end-of-cell marker(s) Word appends to Range.Text\n    $currentText = $cellRange.Text.TrimEnd([char]7, [char]13)\n    if ($currentText -ne $item.OldText) {\n        throw \"Cell ($($item.Row),$($item.Col)) expected `\"$($item.OldText)`\" but found `\"$currentText`\"\"\n    }\n    $cellRange.Text = $item.NewText\n}"}
